$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Lowercase the header row text (E1:O1) in place
$headerRange = $ws.Range("E1:O1")
foreach ($cell in $headerRange.Cells) {
    $cell.Value2 = $cell.Value2.ToLower()
}

# Update active selection to reflect the saved state (F2)
$ws.Range("F2").Select()
